# Update "F" column (想去人数 / "want to go" count) values per commit 456a3b4
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(5, 6).Value = 204
$ws.Cells.Item(7, 6).Value = 120
$ws.Cells.Item(8, 6).Value = 10356
$ws.Cells.Item(10, 6).Value = 3582
$ws.Cells.Item(12, 6).Value = 2463
$ws.Cells.Item(13, 6).Value = 43
$ws.Cells.Item(14, 6).Value = 2863
$ws.Cells.Item(16, 6).Value = 512
$ws.Cells.Item(17, 6).Value = 2198
$ws.Cells.Item(18, 6).Value = 47
$ws.Cells.Item(19, 6).Value = 100
$ws.Cells.Item(24, 6).Value = 322
$ws.Cells.Item(25, 6).Value = 278
$ws.Cells.Item(26, 6).Value = 247
$ws.Cells.Item(27, 6).Value = 620
$ws.Cells.Item(28, 6).Value = 1335
$ws.Cells.Item(29, 6).Value = 19
$ws.Cells.Item(30, 6).Value = 1265
$ws.Cells.Item(32, 6).Value = 133
$ws.Cells.Item(34, 6).Value = 3921
$ws.Cells.Item(35, 6).Value = 3338
$ws.Cells.Item(36, 6).Value = 40
$ws.Cells.Item(38, 6).Value = 1053
$ws.Cells.Item(39, 6).Value = 411
$ws.Cells.Item(42, 6).Value = 114
$ws.Cells.Item(44, 6).Value = 76
$ws.Cells.Item(47, 6).Value = 17

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 11
$ws.Cells.Item(15, 6).Value = 38

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 759
$ws.Cells.Item(3, 6).Value = 1002
$ws.Cells.Item(5, 6).Value = 2098

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(3, 6).Value = 759
$ws.Cells.Item(4, 6).Value = 1002
$ws.Cells.Item(7, 6).Value = 11
$ws.Cells.Item(9, 6).Value = 204
$ws.Cells.Item(11, 6).Value = 120
$ws.Cells.Item(12, 6).Value = 10356
$ws.Cells.Item(15, 6).Value = 3582
$ws.Cells.Item(16, 6).Value = 2463
$ws.Cells.Item(17, 6).Value = 43
$ws.Cells.Item(18, 6).Value = 2863
$ws.Cells.Item(19, 6).Value = 512
$ws.Cells.Item(20, 6).Value = 2198
$ws.Cells.Item(21, 6).Value = 47
$ws.Cells.Item(22, 6).Value = 100
$ws.Cells.Item(26, 6).Value = 322
$ws.Cells.Item(27, 6).Value = 278
$ws.Cells.Item(28, 6).Value = 247
$ws.Cells.Item(29, 6).Value = 620
$ws.Cells.Item(30, 6).Value = 1335
$ws.Cells.Item(31, 6).Value = 19
$ws.Cells.Item(32, 6).Value = 1265
$ws.Cells.Item(34, 6).Value = 133
$ws.Cells.Item(36, 6).Value = 3921
$ws.Cells.Item(37, 6).Value = 3338
$ws.Cells.Item(38, 6).Value = 40
$ws.Cells.Item(39, 6).Value = 1053
$ws.Cells.Item(41, 6).Value = 411
$ws.Cells.Item(43, 6).Value = 38
$ws.Cells.Item(45, 6).Value = 114
$ws.Cells.Item(46, 6).Value = 76
$ws.Cells.Item(48, 6).Value = 17

